# Fixed a bug in sumSymbolVals
# The rows describing each symbol (A3:F23) were being produced in the
# wrong order. This re-writes that block with the corrected row order
# while leaving the header row, the "601" row, and the trailing rows /
# totals row untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(201, 9, 30, 15, 45, 30),
    @(1001, 18, 30, 75, 60, 72),
    @(301, 6, 45, 30, 60, 45),
    @(101, 9, 30, 15, 60, 15),
    @(1203, 3, 15, 15, 15, 15),
    @(701, 3, 90, 45, 97, 15),
    @(801, 3, 67, 65, 52, 45),
    @(901, 16, 15, 45, 60, 60),
    @(401, 9, 48, 67, 75, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(902, 1, 0, 0, 0, 0),
    @(501, 9, 52, 30, 75, 45),
    @(2, 0, 2, 2, 2, 2),
    @(1101, 0, 15, 30, 30, 0),
    @(1, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(802, 0, 4, 5, 4, 0),
    @(502, 0, 4, 0, 0, 0),
    @(402, 0, 0, 4, 0, 0),
    @(602, 0, 0, 4, 0, 9)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
